$d = $word.ActiveDocument

# --- Locate the "No do Laudo: ..." paragraph -------------------------------
# It is the paragraph that carries the Texto1 form-field bookmark
# (<w:bookmarkStart w:name="Texto1"/> ... <w:bookmarkEnd w:id="0"/>).
# We want to insert a brand-new, empty paragraph right after it and
# before the "Tecnico Responsavel:" paragraph, matching:
#
#   <w:p>
#     <w:pPr>
#       <w:spacing w:after="0"/>
#       <w:rPr>
#         <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
#         <w:sz w:val="21"/>
#         <w:szCs w:val="21"/>
#       </w:rPr>
#     </w:pPr>
#   </w:p>
$bk = $d.Bookmarks.Item("Texto1")
$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($bk.Range.Start -ge $p.Range.Start -and $bk.Range.Start -lt $p.Range.End) {
        $targetIndex = $i
        break
    }
}

# --- Find a pre-existing, completely empty paragraph (just the pilcrow) ----
# Its own Range.FormattedText is used as a donor below purely so the
# paragraph-mark run we insert round-trips without Word fabricating a
# stray, content-less <w:r> inside the new paragraph. The donor's own
# formatting is irrelevant -- only its emptiness matters, since the
# destination paragraph already inherits the correct pPr/rPr from the
# InsertParagraphAfter() call (it splits off of the "No do Laudo:" run
# properties, which already use Arial/21).
$donorFormattedText = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text -eq "`r") {
        $donorFormattedText = $cand.Range.FormattedText
        break
    }
}

$target = $d.Paragraphs.Item($targetIndex)
$target.Range.InsertParagraphAfter()

$newPara = $d.Paragraphs.Item($targetIndex + 1)
$newPara.Range.FormattedText = $donorFormattedText

# The FormattedText assignment splits the paragraph mark in two; drop
# the now-redundant extra blank paragraph it leaves behind.
$leftover = $d.Paragraphs.Item($targetIndex + 2)
$leftover.Range.Delete()
